$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.901.37"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.814.36"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'309.31"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.3662"
$ws.Range("D9").Value = "'0.07346"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.830.45"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "'5.382"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'0.07090"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "'6.518"
$ws.Range("D16").Value = "'91.63"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'0.000008705"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'14.66"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "26.940.58"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'5.302"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "2.062.52"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'1.895"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'150.94"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'18.29"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "'2.152"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'5.262"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'115.44"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'0.08919"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'0.7541"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'4.486"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "'2.914"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").Value = "'0.05275"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.01949"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'2.969"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "'7.238"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'0.5302"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'2.277"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").Value = "'0.1653"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "'8.423"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'0.4874"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").Value = "'10.37"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'1.658"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +0.11%  "
